# Adds the new survey wave "25. 1. 2022" as the next column on both worksheets
# ("data" gets new column AN, "pocetR" gets new column AM), matching the
# pattern already used for every previous wave, and refreshes the two
# trailing footnote cells that record the "aktualizace" (last update) date.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "data"
$ws2 = $wb.Worksheets.Item(2)   # "pocetR"

# ----------------------------------------------------------------------
# Sheet "data": new column AN = wave of 25. 1. 2022
# ----------------------------------------------------------------------
$ws1.Range("AN1").Value = "25. 1. 2022"
$ws1.Range("AN2").Value = 0.27
$ws1.Range("AN3").Value = 0.42
$ws1.Range("AN4").Value = 0.31
$ws1.Range("AN5").Value = 0.36
$ws1.Range("AN6").Value = 0.42
$ws1.Range("AN7").Value = 0.22
$ws1.Range("AN8").Value = 0.25
$ws1.Range("AN9").Value = 0.42
$ws1.Range("AN10").Value = 0.33
$ws1.Range("AN11").Value = 0.36
$ws1.Range("AN12").Value = 0.44
$ws1.Range("AN13").Value = 0.2
$ws1.Range("AN14").Value = 0.3
$ws1.Range("AN15").Value = 0.45
$ws1.Range("AN16").Value = 0.25
$ws1.Range("AN17").Value = 0.33
$ws1.Range("AN18").Value = 0.36
$ws1.Range("AN19").Value = 0.31
$ws1.Range("AN20").Value = 0.18
$ws1.Range("AN21").Value = 0.37
$ws1.Range("AN22").Value = 0.45
$ws1.Range("AN23").Value = 0.37
$ws1.Range("AN24").Value = 0.42
$ws1.Range("AN25").Value = 0.21
$ws1.Range("AN26").Value = 0.28
$ws1.Range("AN27").Value = 0.47
$ws1.Range("AN28").Value = 0.25
$ws1.Range("AN29").Value = 0.2
$ws1.Range("AN30").Value = 0.37
$ws1.Range("AN31").Value = 0.43
$ws1.Range("AN32").Value = 0.33
$ws1.Range("AN33").Value = 0.41
$ws1.Range("AN34").Value = 0.26
$ws1.Range("AN35").Value = 0.2
$ws1.Range("AN36").Value = 0.4
$ws1.Range("AN37").Value = 0.4
$ws1.Range("AN38").Value = 0.27
$ws1.Range("AN39").Value = 0.44
$ws1.Range("AN40").Value = 0.29
$ws1.Range("AN41").Value = 0.35
$ws1.Range("AN42").Value = 0.39
$ws1.Range("AN43").Value = 0.26
$ws1.Range("AN44").Value = 0.2
$ws1.Range("AN45").Value = 0.45
$ws1.Range("AN46").Value = 0.35
$ws1.Range("AN47").Value = 0.25
$ws1.Range("AN48").Value = 0.44
$ws1.Range("AN49").Value = 0.31
$ws1.Range("AN50").Value = 0.32
$ws1.Range("AN51").Value = 0.4
$ws1.Range("AN52").Value = 0.28
$ws1.Range("AN53").Value = 0.29
$ws1.Range("AN54").Value = 0.37
$ws1.Range("AN55").Value = 0.34
$ws1.Range("AN56").Value = 0.25
$ws1.Range("AN57").Value = 0.39
$ws1.Range("AN58").Value = 0.36
$ws1.Range("AN59").Value = 0.03
$ws1.Range("AN60").Value = 0.52
$ws1.Range("AN61").Value = 0.45
$ws1.Range("AN62").Value = 0.31
$ws1.Range("AN63").Value = 0.4
$ws1.Range("AN64").Value = 0.29
$ws1.Range("AN65").Value = 0.22
$ws1.Range("AN66").Value = 0.33
$ws1.Range("AN67").Value = 0.45
$ws1.Range("AN68").Value = 0.53
$ws1.Range("AN69").Value = 0.16
$ws1.Range("AN70").Value = 0.31
$ws1.Range("AN71").Value = 0.37
$ws1.Range("AN72").Value = 0.41
$ws1.Range("AN73").Value = 0.22
$ws1.Range("AN74").Value = 0.29
$ws1.Range("AN75").Value = 0.47
$ws1.Range("AN76").Value = 0.24
$ws1.Range("AN77").Value = 0.51
$ws1.Range("AN78").Value = 0.3
$ws1.Range("AN79").Value = 0.19
$ws1.Range("AN80").Value = 0.31
$ws1.Range("AN81").Value = 0.45
$ws1.Range("AN82").Value = 0.24
$ws1.Range("AN83").Value = 0.26
$ws1.Range("AN84").Value = 0.38
$ws1.Range("AN85").Value = 0.36
$ws1.Range("AN86").Value = 0.22
$ws1.Range("AN87").Value = 0.55
$ws1.Range("AN88").Value = 0.23
$ws1.Range("AN89").Value = 0.24
$ws1.Range("AN90").Value = 0.48
$ws1.Range("AN91").Value = 0.28
$ws1.Range("AN92").Value = 0.15
$ws1.Range("AN93").Value = 0.36
$ws1.Range("AN94").Value = 0.49

# Match the formatting of the previous wave column (AM) on the new one (AN)
$ws1.Range("AM1").Copy()
$ws1.Range("AN1").PasteSpecial(-4122)     # xlPasteFormats (header cell)
$ws1.Range("AM2:AM94").Copy()
$ws1.Range("AN2:AN94").PasteSpecial(-4122) # xlPasteFormats (data cells)
$excel.CutCopyMode = 0

# Refresh the footnote with the new "aktualizace" date
$ws1.Range("A95").Value = "Život během pandemie, Obavy z epidemie, % respondentů celkově a ve skupinách, aktualizace 1. 2. 2022"

# ----------------------------------------------------------------------
# Sheet "pocetR": new column AM = wave of 25. 1. 2022
# ----------------------------------------------------------------------
$ws2.Range("AM1").Value = "25. 1. 2022"
$ws2.Range("AM2").Value = 1815
$ws2.Range("AM3").Value = 412
$ws2.Range("AM4").Value = 1403
$ws2.Range("AM5").Value = 304
$ws2.Range("AM6").Value = 797
$ws2.Range("AM7").Value = 109
$ws2.Range("AM8").Value = 605
$ws2.Range("AM9").Value = 444
$ws2.Range("AM10").Value = 676
$ws2.Range("AM11").Value = 695
$ws2.Range("AM12").Value = 666
$ws2.Range("AM13").Value = 500
$ws2.Range("AM14").Value = 649
$ws2.Range("AM15").Value = 880
$ws2.Range("AM16").Value = 935
$ws2.Range("AM17").Value = 952
$ws2.Range("AM18").Value = 413
$ws2.Range("AM19").Value = 211
$ws2.Range("AM20").Value = 239
$ws2.Range("AM21").Value = 41
$ws2.Range("AM22").Value = 153
$ws2.Range("AM23").Value = 84
$ws2.Range("AM24").Value = 13
$ws2.Range("AM25").Value = 254
$ws2.Range("AM26").Value = 486
$ws2.Range("AM27").Value = 232
$ws2.Range("AM28").Value = 338
$ws2.Range("AM29").Value = 310
$ws2.Range("AM30").Value = 212
$ws2.Range("AM31").Value = 338
$ws2.Range("AM32").Value = 385

# Match the formatting of the previous wave column (AL) on the new one (AM)
$ws2.Range("AL1").Copy()
$ws2.Range("AM1").PasteSpecial(-4122)      # xlPasteFormats (header cell)
$ws2.Range("AL2:AL32").Copy()
$ws2.Range("AM2:AM32").PasteSpecial(-4122)  # xlPasteFormats (data cells)
$excel.CutCopyMode = 0

# Refresh the footnote with the new "aktualizace" date (row 33 stays blank in col AM,
# same as every other trailing placeholder cell on this row)
$ws2.Range("A33").Value = "Život během pandemie, Obavy z epidemie, velikost dotázaného souboru celkově a ve skupinách, aktualizace 1. 2. 2022"
